$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Now update the values of rows 133-135 per the new weekly data,
# and set the values for the newly added row 136.

# Row 133: date moves forward
$ws.Range("D133").Value = 44448

# Row 134: date, quality and prices change
$ws.Range("D134").Value = 44167
$ws.Range("I134").Value = "Primera"
$ws.Range("K134").Value = 600
$ws.Range("L134").Value = 600
$ws.Range("M134").Value = 600
$ws.Range("P134").Value = 600

# Row 135: date changes
$ws.Range("D135").Value = 44399

# Row 136 (new row): same data the old row 135 used to hold
$ws.Range("A136").Value = 5
$ws.Range("B136").Value = "Macroferia Regional de Talca"
$ws.Range("C136").Value = "Maule"
$ws.Range("D136").Value = 44400
$ws.Range("D136").NumberFormat = $ws.Range("D135").NumberFormat
$ws.Range("E136").Value = 7
$ws.Range("F136").Value = 100112008
$ws.Range("G136").Value = "Coliflor"
$ws.Range("H136").Value = "Sin especificar"
$ws.Range("I136").Value = "Segunda"
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 500
$ws.Range("L136").Value = 500
$ws.Range("M136").Value = 500
$ws.Range("N136").Value = "`$/unidad"
$ws.Range("O136").Value = "Región del Maule"
$ws.Range("P136").Value = 500
$ws.Range("Q136").Value = 1
$ws.Range("R136").Value = "Hortaliza"
